# Two-digit / one-digit division worksheet — refresh the generated problems.
#
# wdReplaceOne = 1 (replace first match only, leave rest of doc untouched)
# Using a fresh $d.Content range before each call re-scans the whole story,
# so repeated calls walk forward through duplicate matches in document order.

$d = $word.ActiveDocument

function Replace-FirstMatch([string]$old, [string]$new) {
    $rng = $d.Content
    $rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 1) | Out-Null
}

# Row 1
Replace-FirstMatch "76÷9=" "86÷8="
Replace-FirstMatch "34÷2=" "69÷6="
Replace-FirstMatch "54÷3=" "34÷9="
Replace-FirstMatch "53÷8=" "99÷6="
Replace-FirstMatch "39÷4=" "34÷7="

# Row 2
Replace-FirstMatch "78÷9=" "52÷6="
Replace-FirstMatch "85÷6=" "27÷5="
Replace-FirstMatch "42÷7=" "71÷7="
Replace-FirstMatch "66÷4=" "89÷6="
# second (and last) remaining occurrence of "39÷4="
Replace-FirstMatch "39÷4=" "90÷2="

# Row 3
Replace-FirstMatch "59÷3=" "19÷5="
Replace-FirstMatch "12÷4=" "99÷9="
Replace-FirstMatch "30÷6=" "20÷5="
Replace-FirstMatch "50÷4=" "32÷2="
Replace-FirstMatch "24÷8=" "52÷2="

# Row 4 — the "10÷9=" cell is dropped and a new "88÷9=" cell is appended at
# the end of the row, so cells 3-5 shift one slot to the left before the new
# value lands in the now-last cell. Cell count stays 5, so this is done as a
# straight per-cell text assignment rather than a structural insert/delete.
Replace-FirstMatch "65÷3=" "34÷8="
Replace-FirstMatch "64÷9=" "86÷6="
$table = $d.Tables.Item(1)
$row4 = $table.Rows.Item(13)
$row4.Cells.Item(3).Range.Text = "77÷5="
$row4.Cells.Item(4).Range.Text = "75÷9="
$row4.Cells.Item(5).Range.Text = "88÷9="

# Row 5
Replace-FirstMatch "21÷5=" "40÷3="
Replace-FirstMatch "69÷9=" "79÷8="
Replace-FirstMatch "95÷3=" "17÷3="
Replace-FirstMatch "18÷5=" "81÷9="
Replace-FirstMatch "52÷8=" "68÷4="
